# Extend the TestCase schema workbook:
#  - rename "TestCaseCollection" -> "SemanticSmokeTestCase"
#  - add new sheets: Input, SemanticSmokeTestInput, Output,
#    SemanticSmokeTestOutput, Precondition, TestSuite
#  - add inputs/outputs/preconditions columns to TestCase & SemanticSmokeTestCase
#  - populate SemanticSmokeTestInput headers + its 3 list data validations
#  - TestSuite keeps the old "entries" header (what TestCaseCollection used to hold)

$wb = $excel.ActiveWorkbook

# --- rename sheet2 (TestCaseCollection -> SemanticSmokeTestCase) ---
$testCaseCollection = $wb.Worksheets.Item("TestCaseCollection")
$testCaseCollection.Name = "SemanticSmokeTestCase"

# --- add the rest of the sheets, in tab order, after the current last sheet ---
$sheetNames = @("Input", "SemanticSmokeTestInput", "Output", "SemanticSmokeTestOutput", "Precondition", "TestSuite")
foreach ($name in $sheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
}

# --- TestCase: extend header row with inputs / outputs / preconditions ---
$testCase = $wb.Worksheets.Item("TestCase")
$testCase.Range("D1").Value = "inputs"
$testCase.Range("E1").Value = "outputs"
$testCase.Range("F1").Value = "preconditions"

# --- SemanticSmokeTestCase: same header layout as TestCase ---
$semanticSmokeTestCase = $wb.Worksheets.Item("SemanticSmokeTestCase")
$semanticSmokeTestCase.Range("A1").Value = "name"
$semanticSmokeTestCase.Range("B1").Value = "id"
$semanticSmokeTestCase.Range("C1").Value = "description"
$semanticSmokeTestCase.Range("D1").Value = "inputs"
$semanticSmokeTestCase.Range("E1").Value = "outputs"
$semanticSmokeTestCase.Range("F1").Value = "preconditions"

# --- SemanticSmokeTestInput: headers + data validations ---
$semanticSmokeTestInput = $wb.Worksheets.Item("SemanticSmokeTestInput")
$semanticSmokeTestInput.Range("A1").Value = "must_pass_date"
$semanticSmokeTestInput.Range("B1").Value = "must_pass_environment"
$semanticSmokeTestInput.Range("C1").Value = "query"
$semanticSmokeTestInput.Range("D1").Value = "string_entry"
$semanticSmokeTestInput.Range("E1").Value = "direction"
$semanticSmokeTestInput.Range("F1").Value = "answer_informal_concept"
$semanticSmokeTestInput.Range("G1").Value = "expected_result"
$semanticSmokeTestInput.Range("H1").Value = "curie"
$semanticSmokeTestInput.Range("I1").Value = "top_level"
$semanticSmokeTestInput.Range("J1").Value = "node"
$semanticSmokeTestInput.Range("K1").Value = "notes"

$semanticSmokeTestInput.Range("B2:B1048576").Validation.Add(3, 1, 1, '"DEV,CI,TEST,PROD"')
$semanticSmokeTestInput.Range("E2:E1048576").Validation.Add(3, 1, 1, '"increased,decreased"')
$semanticSmokeTestInput.Range("G2:G1048576").Validation.Add(3, 1, 1, '"include_good,exclude_bad"')

# --- TestSuite: carries over the old TestCaseCollection content ---
$testSuite = $wb.Worksheets.Item("TestSuite")
$testSuite.Range("A1").Value = "entries"
